$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "872.0 (11)"
$ws.Range("E4").Value = "750.5 (11)"
$ws.Range("F4").Value = "121.5 (0)"

$ws.Range("D5").Value = "857.0 (27)"
$ws.Range("E5").Value = "744.5 (18)"
$ws.Range("F5").Value = "112.5 (9)"

$ws.Range("D6").Value = "587.0 (36)"
$ws.Range("E6").Value = "546.5 (22)"
$ws.Range("F6").Value = "40.5 (14)"

$ws.Range("D7").Value = "1084.0 (15)"
$ws.Range("E7").Value = "909.5 (15)"
$ws.Range("F7").Value = "174.5 (0)"

$ws.Range("D8").Value = "1282.0 (33)"
$ws.Range("E8").Value = "998.0 (28)"
$ws.Range("F8").Value = "284.0 (5)"

$ws.Range("D9").Value = "1355.0 (32)"
$ws.Range("E9").Value = "909.0 (23)"
$ws.Range("F9").Value = "446.0 (8)"

$ws.Range("D10").Value = "931.0 (18)"
$ws.Range("E10").Value = "829.0 (11)"
$ws.Range("F10").Value = "102.0 (7)"

$ws.Range("D11").Value = "838.0 (36)"
$ws.Range("E11").Value = "772.5 (21)"
$ws.Range("F11").Value = "65.5 (15)"

$ws.Range("D12").Value = "792.0 (64)"
$ws.Range("E12").Value = "650.5 (35)"
$ws.Range("F12").Value = "141.5 (29)"
